# Update "想去人数" (want-to-go count) values in column F across sheets
# to reflect the refreshed data output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1151
$ws1.Range("F10").Value = 67
$ws1.Range("F12").Value = 413
$ws1.Range("F18").Value = 6249
$ws1.Range("F27").Value = 798
$ws1.Range("F28").Value = 4477
$ws1.Range("F32").Value = 1292
$ws1.Range("F33").Value = 122
$ws1.Range("F36").Value = 1008
$ws1.Range("F37").Value = 1316

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1170

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1170
$ws4.Range("F10").Value = 1151
$ws4.Range("F13").Value = 67
$ws4.Range("F15").Value = 413
$ws4.Range("F22").Value = 6249
$ws4.Range("F23").Value = 6249
$ws4.Range("F32").Value = 798
$ws4.Range("F33").Value = 4477
$ws4.Range("F38").Value = 1292
$ws4.Range("F39").Value = 122
$ws4.Range("F42").Value = 1008
$ws4.Range("F43").Value = 1316
